$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.466.26"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.668.43"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").Value = "'313.56"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.3963"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "'52.01"
$ws.Range("E9").Value = "  +6.48%  "
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("D11").Value = "'0.9993"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "'0.08607"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").Value = "'24.53"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "'7.362"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "'0.00001348"
$ws.Range("E15").Value = "  +4.78%  "
$ws.Range("D16").Value = "'7.924"
$ws.Range("E16").Value = "  +5.34%  "
$ws.Range("D17").Value = "1.670.10"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'95.60"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "'0.06985"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "'20.66"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'13.79"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "24.454.55"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'2.428"
$ws.Range("E25").Value = "  +3.67%  "
$ws.Range("D26").Value = "'3.037"
$ws.Range("E26").Value = "  +11.69%  "
$ws.Range("D27").Value = "'22.55"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'157.66"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "'143.07"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").Value = "'5.473"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").Value = "'8.161"
$ws.Range("E31").Value = "  -8.93%  "
$ws.Range("D32").Value = "'2.547"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "1.849.19"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'1.063"
$ws.Range("E34").Value = "  +7.40%  "
$ws.Range("D35").Value = "'0.08292"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").Value = "'0.03031"
$ws.Range("E36").Value = "  +2.68%  "
$ws.Range("D37").Value = "'6.848"
$ws.Range("E37").Value = "  -5.46%  "
$ws.Range("D38").Value = "'11.16"
$ws.Range("E38").Value = "  +10.93%  "
$ws.Range("D39").Value = "'0.2770"
$ws.Range("D40").Value = "'0.09271"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").Value = "'0.7774"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "'13.89"
$ws.Range("E42").Value = "  +6.07%  "
$ws.Range("D43").Value = "'1.444"
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").Value = "'16.67"
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("D45").Value = "'0.7152"
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").Value = "'2.543"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "'4.148"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Value = "'0.9996"
$ws.Range("D49").Value = "'0.08460"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'136.87"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "'1.285"
$ws.Range("E51").Value = "  +1.70%  "
